# ---------------------------------------------------------------------------
# "Generate Report for Handback"
#
# The CI job re-ran and produced a fresh handback report: two new source
# files (identified by new GUID-based names) replaced the two files that
# used to be tracked, the handoff/handback xliff artifacts were
# regenerated with a new content hash, and all of the "generated at"
# timestamps were bumped forward by about a minute. This script pushes
# those new values into the three report worksheets (Overview, zh-cn,
# de-de) and refreshes the corresponding hyperlinks.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldGuid1 = "2c5e234e-fb30-4425-885b-b108c390ff7c"
$newGuid1 = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0"
$oldGuid2 = "abf2a8c3-5959-415f-906b-fab2aa37853a"
$newGuid2 = "ffff341494fd-1949-4346-98a4-80679563b85d"
$newHash  = "5699311b2f211a1627904973773dd11918028c8f"

# ---------------------------------------------------------------------------
# Sheet "Overview": update file name / path / "Latest HO Xliff Generate
# Date" cells for the two handback rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-09-05 23:16:38"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-09-05 23:16:38"

# Refresh the hyperlinks backing column B (this engine only lets us clear
# *all* hyperlinks on the sheet at once, so drop them all and re-add the
# two we need with their updated display text / target URL).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": update file name / handoff-handback xliff file name and
# datetime cells for the two handback rows.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-05 23:16:32"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 23:16:50"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 23:16:32"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 23:16:50"

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5f60fdfb60bc7385f0fb3c0376e9514ac192e11/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5f60fdfb60bc7385f0fb3c0376e9514ac192e11/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": update file name / handoff-handback xliff file name and
# datetime cells for the two handback rows.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-05 23:16:38"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 23:16:58"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 23:16:38"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 23:16:58"

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4d66e88fe7f9e24d6a91b0ca29def2da65c52cbe/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4d66e88fe7f9e24d6a91b0ca29def2da65c52cbe/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null
